# Auto-generated Excel COM-interop script to add new data rows
# matching the target OOXML diff for content-data-template.xlsx

$wb = $excel.ActiveWorkbook

# ----- Sheet: FinancialRecord -----
$ws = $wb.Worksheets.Item('FinancialRecord')

# Row 4
$ws.Cells.Item(4, 1).Value = 'cuid_fin_202403'  # A4
$c = $ws.Cells.Item(4, 2)  # B4 (date-like text; force text to avoid auto date conversion)
$c.NumberFormat = '@'
$c.Value = '2024-03-01'
$ws.Cells.Item(4, 3).Value = 125000.5  # C4
$ws.Cells.Item(4, 4).Value = 97250.75  # D4
$ws.Cells.Item(4, 5).Value = 27749.75  # E4
$ws.Cells.Item(4, 6).Value = 'Mar 2024 Ministry Update'  # F4
$ws.Cells.Item(4, 7).Value = '[{"id":"tithes","label":"ถวายสิบลด","amount":82000},{"id":"partnerships","label":"Partnership Gifts","amount":43000}]'  # G4
$ws.Cells.Item(4, 8).Value = '[{"id":"operations","label":"ค่าใช้จ่ายดำเนินงาน","amount":32000},{"id":"missions","label":"Missions Support","amount":22000},{"id":"staff","label":"Staff Salaries","amount":43250.75}]'  # H4

# Row 5
$c = $ws.Cells.Item(5, 2)  # B5 (date-like text; force text to avoid auto date conversion)
$c.NumberFormat = '@'
$c.Value = '2024-04-01'
$ws.Cells.Item(5, 3).Value = 118500  # C5
$ws.Cells.Item(5, 4).Value = 110200  # D5
$ws.Cells.Item(5, 5).Value = 20049.75  # E5
$ws.Cells.Item(5, 6).Value = 'Apr 2024 Easter Activities'  # F5
$ws.Cells.Item(5, 7).Value = '[{"id":"donations","label":"Easter Offering","amount":45500},{"id":"grants","label":"Community Grant","amount":73000}]'  # G5
$ws.Cells.Item(5, 8).Value = '[{"id":"outreach","label":"Outreach Events","amount":42000},{"id":"benevolence","label":"Benevolence","amount":18500},{"id":"facilities","label":"Facility Upgrades","amount":49600}]'  # H5


# ----- Sheet: FutureProject -----
$ws = $wb.Worksheets.Item('FutureProject')

# Row 4
$ws.Cells.Item(4, 1).Value = 'cuid_proj_center'  # A4
$ws.Cells.Item(4, 2).Value = 'Community Center Renovation'  # B4
$ws.Cells.Item(4, 3).Value = 'Upgrade the main hall and classrooms to expand youth programs.'  # C4
$ws.Cells.Item(4, 4).Value = 500000  # D4
$ws.Cells.Item(4, 5).Value = 185000  # E4
$ws.Cells.Item(4, 6).Value = 1  # F4
$ws.Cells.Item(4, 7).Value = $true  # G4

# Row 5
$ws.Cells.Item(5, 2).Value = 'Mobile Medical Clinic'  # B5
$ws.Cells.Item(5, 3).Value = 'Equip a mobile team to provide basic care in rural villages.'  # C5
$ws.Cells.Item(5, 4).Value = 350000  # D5
$ws.Cells.Item(5, 5).Value = 45000  # E5
$ws.Cells.Item(5, 6).Value = 2  # F5
$ws.Cells.Item(5, 7).Value = $true  # G5


# ----- Sheet: Mission -----
$ws = $wb.Worksheets.Item('Mission')

# Row 4
$ws.Cells.Item(4, 1).Value = 'cuid_mission_youth'  # A4
$ws.Cells.Item(4, 2).Value = 'chiang-mai-youth'  # B4
$ws.Cells.Item(4, 3).Value = 'พันธกิจเยาวชนเชียงใหม่'  # C4
$ws.Cells.Item(4, 4).Value = 'Chiang Mai Youth Outreach'  # D4
$ws.Cells.Item(4, 5).Value = 'สร้างผู้นำรุ่นใหม่'  # E4
$ws.Cells.Item(4, 6).Value = 'Raising Young Leaders'  # F4
$ws.Cells.Item(4, 7).Value = 'ทำงานร่วมกับคริสตจักรท้องถิ่นเพื่อเสริมสร้างเยาวชนให้เติบโตเป็นสาวก'  # G4
$ws.Cells.Item(4, 8).Value = 'Partner with local churches to equip students as disciples.'  # H4
$ws.Cells.Item(4, 9).Value = 'ทีมงานจัดค่าย การอบรม และการติดตามผลเพื่อสนับสนุนเยาวชนและครอบครัวตลอดปี'  # I4
$ws.Cells.Item(4, 10).Value = 'The team runs camps, leadership labs, and follow-up coaching with families year-round.'  # J4
$ws.Cells.Item(4, 11).Value = 'ค่ายเยาวชน | การเป็นผู้นำ | ศิลปะสร้างสรรค์'  # K4
$ws.Cells.Item(4, 12).Value = 'Youth Camps | Leadership Labs | Creative Arts'  # L4
$ws.Cells.Item(4, 13).Value = 'มัทธิว 5:14'  # M4
$ws.Cells.Item(4, 14).Value = 'Matthew 5:14'  # N4
$ws.Cells.Item(4, 15).Value = 'ท่านทั้งหลายเป็นความสว่างของโลก...'  # O4
$ws.Cells.Item(4, 16).Value = 'You are the light of the world...'  # P4
$ws.Cells.Item(4, 17).Value = 'อธิษฐานสนับสนุน | ร่วมอาสา | ให้การสนับสนุนทางการเงิน'  # Q4
$ws.Cells.Item(4, 18).Value = 'Pray with us | Volunteer on-site | Become a monthly partner'  # R4
$ws.Cells.Item(4, 19).Value = $true  # S4
$ws.Cells.Item(4, 20).Value = 'https://example.org/images/chiang-mai-youth.jpg'  # T4
$c = $ws.Cells.Item(4, 21)  # U4 (date-like text; force text to avoid auto date conversion)
$c.NumberFormat = '@'
$c.Value = '2024-01-15'
$c = $ws.Cells.Item(4, 22)  # V4 (date-like text; force text to avoid auto date conversion)
$c.NumberFormat = '@'
$c.Value = '2024-12-15'


# ----- Sheet: ContactInfo -----
$ws = $wb.Worksheets.Item('ContactInfo')

# Row 4
$ws.Cells.Item(4, 1).Value = 1  # A4
$ws.Cells.Item(4, 2).Value = 'คริสตจักรความหวังกรุงเทพ'  # B4
$ws.Cells.Item(4, 3).Value = 'Hope Church Bangkok'  # C4
$ws.Cells.Item(4, 4).Value = '02-123-4567, 081-234-5678'  # D4
$ws.Cells.Item(4, 5).Value = 'office@hopebkk.org'  # E4
$ws.Cells.Item(4, 6).Value = '123 ถนนสุขุมวิท แขวงคลองตัน เขตคลองเตย กรุงเทพฯ 10110'  # F4
$ws.Cells.Item(4, 7).Value = '123 Sukhumvit Rd, Khlong Toei, Bangkok 10110'  # G4
$ws.Cells.Item(4, 8).Value = 'https://facebook.com/hopechurchbkk'  # H4
$ws.Cells.Item(4, 9).Value = 'https://facebook.com/hopechurchbkk/live'  # I4
$ws.Cells.Item(4, 10).Value = 'https://youtube.com/@hopechurchbkk'  # J4
$ws.Cells.Item(4, 11).Value = 'https://maps.google.com/?q=13.73,100.567'  # K4
$ws.Cells.Item(4, 12).Value = 13.73  # L4
$ws.Cells.Item(4, 13).Value = 100.567  # M4
$ws.Cells.Item(4, 14).Value = "Sunday|09:30|Thai Service`nSunday|11:30|English Service`nWednesday|19:00|Prayer Gathering"  # N4


# ----- Sheet: NavigationItem -----
$ws = $wb.Worksheets.Item('NavigationItem')

# Row 4
$ws.Cells.Item(4, 1).Value = 'cuid_nav_home'  # A4
$ws.Cells.Item(4, 2).Value = 'หน้าแรก'  # B4
$ws.Cells.Item(4, 3).Value = 'Home'  # C4
$ws.Cells.Item(4, 4).Value = '/'  # D4
$ws.Cells.Item(4, 5).Value = 1  # E4
$ws.Cells.Item(4, 6).Value = $true  # F4

# Row 5
$ws.Cells.Item(5, 1).Value = 'cuid_nav_missions'  # A5
$ws.Cells.Item(5, 2).Value = 'พันธกิจ'  # B5
$ws.Cells.Item(5, 3).Value = 'Missions'  # C5
$ws.Cells.Item(5, 4).Value = '/missions'  # D5
$ws.Cells.Item(5, 5).Value = 2  # E5
$ws.Cells.Item(5, 6).Value = $true  # F5

# Row 6
$ws.Cells.Item(6, 2).Value = 'ถวาย'  # B6
$ws.Cells.Item(6, 3).Value = 'Give'  # C6
$ws.Cells.Item(6, 4).Value = '/give'  # D6
$ws.Cells.Item(6, 5).Value = 3  # E6
$ws.Cells.Item(6, 6).Value = $true  # F6


# ----- Sheet: PageContent -----
$ws = $wb.Worksheets.Item('PageContent')

# Row 4
$ws.Cells.Item(4, 1).Value = 'cuid_page_landing_hero'  # A4
$ws.Cells.Item(4, 2).Value = 'landing'  # B4
$ws.Cells.Item(4, 3).Value = 'hero'  # C4
$ws.Cells.Item(4, 4).Value = 'ยินดีต้อนรับสู่คริสตจักร'  # D4
$ws.Cells.Item(4, 5).Value = 'Welcome to Our Church'  # E4
$ws.Cells.Item(4, 6).Value = 'ครอบครัวที่เติบโตไปด้วยกัน'  # F4
$ws.Cells.Item(4, 7).Value = 'A family growing together'  # G4
$ws.Cells.Item(4, 8).Value = 'เรามุ่งมั่นที่จะเห็นชีวิตได้รับการเปลี่ยนแปลงผ่านการติดตามพระเยซู'  # H4
$ws.Cells.Item(4, 9).Value = 'We pursue transformed lives through following Jesus together.'  # I4
$ws.Cells.Item(4, 10).Value = '{"ctaLabel":"Join Us","ctaHref":"/contact"}'  # J4
$ws.Cells.Item(4, 11).Value = '{"background":"sunrise","metaDescription":"Hope Church Bangkok landing hero"}'  # K4

# Row 5
$ws.Cells.Item(5, 2).Value = 'landing'  # B5
$ws.Cells.Item(5, 3).Value = 'stories'  # C5
$ws.Cells.Item(5, 4).Value = 'คำพยานล่าสุด'  # D5
$ws.Cells.Item(5, 5).Value = 'Latest Stories'  # E5
$ws.Cells.Item(5, 6).Value = 'พระเจ้าทรงทำสิ่งใหม่'  # F5
$ws.Cells.Item(5, 7).Value = 'God is doing something new'  # G5
$ws.Cells.Item(5, 8).Value = 'สมาชิกของเรามีประสบการณ์ที่พระเจ้าทรงเปลี่ยนแปลงชีวิตทุกสัปดาห์'  # H5
$ws.Cells.Item(5, 9).Value = 'Every week our members share how God is changing lives.'  # I5
$ws.Cells.Item(5, 10).Value = '[{"name":"Nok","highlight":"Found freedom from anxiety"},{"name":"Ben","highlight":"Started discipling teens"}]'  # J5
$ws.Cells.Item(5, 11).Value = '{"showFilter":true}'  # K5


# ----- Sheet: FinancialCategory -----
$ws = $wb.Worksheets.Item('FinancialCategory')

# Row 4
$ws.Cells.Item(4, 1).Value = 'cuid_cat_tithes'  # A4
$ws.Cells.Item(4, 2).Value = 'tithes'  # B4
$ws.Cells.Item(4, 3).Value = 'Tithes & Offerings'  # C4
$ws.Cells.Item(4, 4).Value = 'income'  # D4
$ws.Cells.Item(4, 5).Value = $true  # E4
$ws.Cells.Item(4, 7).Value = 1  # G4

# Row 5
$ws.Cells.Item(5, 1).Value = 'cuid_cat_missions'  # A5
$ws.Cells.Item(5, 2).Value = 'missions'  # B5
$ws.Cells.Item(5, 3).Value = 'Missions Support'  # C5
$ws.Cells.Item(5, 4).Value = 'expense'  # D5
$ws.Cells.Item(5, 5).Value = $true  # E5
$ws.Cells.Item(5, 7).Value = 2  # G5

# Row 6
$ws.Cells.Item(6, 2).Value = 'administration'  # B6
$ws.Cells.Item(6, 3).Value = 'Administration'  # C6
$ws.Cells.Item(6, 4).Value = 'expense'  # D6
$ws.Cells.Item(6, 5).Value = $false  # E6
$ws.Cells.Item(6, 6).Value = 'operations'  # F6
$ws.Cells.Item(6, 7).Value = 3  # G6
$ws.Cells.Item(6, 8).Value = 2024  # H6


# ----- Sheet: CategorySettings -----
$ws = $wb.Worksheets.Item('CategorySettings')

# Row 4
$ws.Cells.Item(4, 1).Value = 'cuid_settings_2024'  # A4
$ws.Cells.Item(4, 2).Value = 2024  # B4
$ws.Cells.Item(4, 3).Value = '{"featuredCategories":["tithes","missions"],"monthlyGoal":450000,"showCumulative":true}'  # C4


Write-Output "Added new data rows to all 8 sheets."
